$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 (Arabia Saudita): updated case numbers ---
$ws.Range("B25").Value = 18811
$ws.Range("C25").Value = 1289
$ws.Range("D25").Value = 2531
$ws.Range("E25").Value = 16136
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 144

# --- Row 62: updated case numbers ---
$ws.Range("B62").Value = 2835
$ws.Range("C62").Value = 118
$ws.Range("D62").Value = 720
$ws.Range("E62").Value = 2090

# --- Row 81: updated case numbers ---
$ws.Range("B81").Value = 1399
$ws.Range("C81").Value = 13
$ws.Range("D81").Value = 553
$ws.Range("E81").Value = 781
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 65

# --- Rows 95-101: "Senegal" inserted alphabetically before "Albania",
#     pushing Albania..Kirguistan down by one row with their new case data ---
$ws.Range("A95").Value = "Senegal"
$ws.Range("B95").Value = 736
$ws.Range("C95").Value = 65
$ws.Range("D95").Value = 284
$ws.Range("E95").Value = 443
$ws.Range("F95").Value = 1
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 9

$ws.Range("A96").Value = "Albania"
$ws.Range("B96").Value = 726
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 410
$ws.Range("E96").Value = 288
$ws.Range("F96").Value = 4
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 28

$ws.Range("A97").Value = "Crucero"
$ws.Range("B97").Value = 712
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 645
$ws.Range("E97").Value = 54
$ws.Range("F97").Value = 4
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 13

$ws.Range("A98").Value = "Libano"
$ws.Range("B98").Value = 710
$ws.Range("C98").Value = 3
$ws.Range("D98").Value = 145
$ws.Range("E98").Value = 541
$ws.Range("F98").Value = 44
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 24

$ws.Range("A99").Value = "Niger"
$ws.Range("B99").Value = 696
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 350
$ws.Range("E99").Value = 317
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 29

$ws.Range("A100").Value = "Costa Rica"
$ws.Range("B100").Value = 695
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 264
$ws.Range("E100").Value = 425
$ws.Range("F100").Value = 8
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 6

$ws.Range("A101").Value = "Kirguistan"
$ws.Range("B101").Value = 695
$ws.Range("C101").Value = 13
$ws.Range("D101").Value = 395
$ws.Range("E101").Value = 292
$ws.Range("F101").Value = 13
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 8

# --- Rows 146-147: "Cabo Verde" moved ahead of "Bermudas" ---
$ws.Range("A146").Value = "Cabo Verde"
$ws.Range("B146").Value = 109
$ws.Range("C146").Value = 3
$ws.Range("D146").Value = 1
$ws.Range("E146").Value = 107
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 1

$ws.Range("A147").Value = "Bermudas"
$ws.Range("B147").Value = 109
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 39
$ws.Range("E147").Value = 64
$ws.Range("F147").Value = 10
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 6

# --- Row 175: updated case numbers ---
$ws.Range("B175").Value = 36
$ws.Range("C175").Value = 2
$ws.Range("E175").Value = 29

# --- Rows 184-185: "Granada" moved ahead of "Belice" ---
$ws.Range("A184").Value = "Granada"
$ws.Range("D184").Value = 7
$ws.Range("F184").Value = 4
$ws.Range("H184").Value = 0

$ws.Range("A185").Value = "Belice"
$ws.Range("D185").Value = 5
$ws.Range("F185").Value = 1
$ws.Range("H185").Value = 2

# --- Rows 198-199: "Burundi" moved ahead of "Islas Turcas y Caicos" ---
# (case data for both rows is identical, so only the names swap)
$ws.Range("A198").Value = "Burundi"
$ws.Range("A199").Value = "Islas Turcas y Caicos"
